$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.29347
$ws.Range("H2").Value = 0.8804099999999999
$ws.Range("I2").Value = 0.1501202107524681
$ws.Range("J2").Value = 0.1501202107524681
$ws.Range("M2").Value = 0.655792
$ws.Range("N2").Value = 1.967376
$ws.Range("O2").Value = 0.01246532615150124
$ws.Range("P2").Value = 0.01246532615150124
$ws.Range("Q2").Value = 0.19245527824
$ws.Range("R2").Value = 1.73209750416
$ws.Range("S2").Value = 0.001871297388961619
$ws.Range("T2").Value = 0.001871297388961619
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.29347
$ws.Range("H3").Value = 0.8804099999999999
$ws.Range("I3").Value = 0.1501202107524681
$ws.Range("J3").Value = 0.1501202107524681
$ws.Range("O3").Value = 0.5315769812025607
$ws.Range("P3").Value = 0.5315769812025607
$ws.Range("Q3").Value = 8.207149542653331
$ws.Range("R3").Value = 73.86434588387999
$ws.Range("S3").Value = 0.07980044844928921
$ws.Range("T3").Value = 0.07980044844928921
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.29347
$ws.Range("H4").Value = 0.8804099999999999
$ws.Range("I4").Value = 0.1501202107524681
$ws.Range("J4").Value = 0.1501202107524681
$ws.Range("M4").Value = 23.987612
$ws.Range("N4").Value = 71.962836
$ws.Range("O4").Value = 0.4559576926459381
$ws.Range("P4").Value = 0.4559576926459381
$ws.Range("Q4").Value = 7.039644493639998
$ws.Range("R4").Value = 63.35680044275999
$ws.Range("S4").Value = 0.06844846491421731
$ws.Range("T4").Value = 0.06844846491421733
$ws.Range("I5").Value = 0.2180391153852712
$ws.Range("J5").Value = 0.2180391153852712
$ws.Range("M5").Value = 0.655792
$ws.Range("N5").Value = 1.967376
$ws.Range("O5").Value = 0.01246532615150124
$ws.Range("P5").Value = 0.01246532615150124
$ws.Range("Q5").Value = 0.2795278424426667
$ws.Range("R5").Value = 2.515750581984
$ws.Range("S5").Value = 0.002717928687062218
$ws.Range("T5").Value = 0.002717928687062218
$ws.Range("I6").Value = 0.2180391153852712
$ws.Range("J6").Value = 0.2180391153852712
$ws.Range("O6").Value = 0.5315769812025607
$ws.Range("P6").Value = 0.5315769812025607
$ws.Range("S6").Value = 0.1159045747405793
$ws.Range("T6").Value = 0.1159045747405793
$ws.Range("I7").Value = 0.2180391153852712
$ws.Range("J7").Value = 0.2180391153852712
$ws.Range("M7").Value = 23.987612
$ws.Range("N7").Value = 71.962836
$ws.Range("O7").Value = 0.4559576926459381
$ws.Range("P7").Value = 0.4559576926459381
$ws.Range("Q7").Value = 10.22459168106933
$ws.Range("R7").Value = 92.021325129624
$ws.Range("S7").Value = 0.09941661195762971
$ws.Range("T7").Value = 0.09941661195762973
$ws.Range("G8").Value = 1.235185333333333
$ws.Range("H8").Value = 3.705556
$ws.Range("I8").Value = 0.6318406738622607
$ws.Range("J8").Value = 0.6318406738622606
$ws.Range("M8").Value = 0.655792
$ws.Range("N8").Value = 1.967376
$ws.Range("O8").Value = 0.01246532615150124
$ws.Range("P8").Value = 0.01246532615150124
$ws.Range("Q8").Value = 0.8100246601173334
$ws.Range("R8").Value = 7.290221941056
$ws.Range("S8").Value = 0.007876100075477408
$ws.Range("T8").Value = 0.007876100075477406
$ws.Range("G9").Value = 1.235185333333333
$ws.Range("H9").Value = 3.705556
$ws.Range("I9").Value = 0.6318406738622607
$ws.Range("J9").Value = 0.6318406738622606
$ws.Range("O9").Value = 0.5315769812025607
$ws.Range("P9").Value = 0.5315769812025607
$ws.Range("Q9").Value = 34.54305633815645
$ws.Range("R9").Value = 310.887507043408
$ws.Range("S9").Value = 0.3358719580126923
$ws.Range("T9").Value = 0.3358719580126922
$ws.Range("G10").Value = 1.235185333333333
$ws.Range("H10").Value = 3.705556
$ws.Range("I10").Value = 0.6318406738622607
$ws.Range("J10").Value = 0.6318406738622606
$ws.Range("M10").Value = 23.987612
$ws.Range("N10").Value = 71.962836
$ws.Range("O10").Value = 0.4559576926459381
$ws.Range("P10").Value = 0.4559576926459381
$ws.Range("Q10").Value = 29.62914652409066
$ws.Range("R10").Value = 266.662318716816
$ws.Range("S10").Value = 0.2880926157740911
$ws.Range("T10").Value = 0.2880926157740911
